$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to make edits, then re-protect afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure footer (A59).
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-19 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns with the latest
# model snapshot figures for rows 2-56.
$ws.Range("D2").Value = 0.0143788811454294
$ws.Range("E2").Value = -0.002416313059877551
$ws.Range("D3").Value = 0.05072355146588409
$ws.Range("E3").Value = -0.000148501986214189
$ws.Range("D4").Value = 0.01426821545899833
$ws.Range("E4").Value = 0.005147269087789397
$ws.Range("D5").Value = 0.009759671540727483
$ws.Range("E5").Value = 0.00308721378955501
$ws.Range("D6").Value = 0.01552753196562346
$ws.Range("E6").Value = -0.0005497903924129854
$ws.Range("D7").Value = 0.02010889051591928
$ws.Range("E7").Value = 0.002053990610328515
$ws.Range("D8").Value = 0.004804604445810346
$ws.Range("E8").Value = -0.02967671132653538
$ws.Range("D9").Value = 0.006827878261969972
$ws.Range("E9").Value = -0.02399470461691222
$ws.Range("D10").Value = 0.01449017454420551
$ws.Range("E10").Value = -0.003041054232133633
$ws.Range("D11").Value = 0.008344807915000719
$ws.Range("E11").Value = -0.002707988566270414
$ws.Range("D12").Value = 0.0158221173691521
$ws.Range("E12").Value = -0.007617233991906835
$ws.Range("D13").Value = 0.003075790490710359
$ws.Range("E13").Value = -0.01142857142857145
$ws.Range("D14").Value = 0.005989945052572164
$ws.Range("E14").Value = -0.01959654178674353
$ws.Range("D15").Value = 0.01477681938656936
$ws.Range("E15").Value = -0.00763781952571585
$ws.Range("D16").Value = 0.0108214469723676
$ws.Range("E16").Value = -0.002088227616810268
$ws.Range("D17").Value = 0.02126155869512162
$ws.Range("E17").Value = -0.009233426331871852
$ws.Range("D18").Value = 0.008647365265311409
$ws.Range("E18").Value = -0.0007621951219511924
$ws.Range("D19").Value = 0.0171189710740108
$ws.Range("E19").Value = -0.002170724552654546
$ws.Range("D20").Value = 0.01221716537157055
$ws.Range("E20").Value = 0.009659353645378399
$ws.Range("D21").Value = 0.006863092924526881
$ws.Range("E21").Value = -0.01996615905245358
$ws.Range("D22").Value = 0.01557404545039114
$ws.Range("E22").Value = 0.00677124478049862
$ws.Range("D23").Value = 0.01971842205168533
$ws.Range("E23").Value = -0.001976879109544916
$ws.Range("D24").Value = 0.01031237226053778
$ws.Range("E24").Value = -0.02586054721977049
$ws.Range("D25").Value = 0.0200154869189769
$ws.Range("E25").Value = -0.004240052185257692
$ws.Range("D26").Value = 0.01410403729515235
$ws.Range("E26").Value = 0.001682321603813097
$ws.Range("D27").Value = 0.02019077559133014
$ws.Range("E27").Value = 0.02970093873433322
$ws.Range("D28").Value = 0.05485892039502691
$ws.Range("E28").Value = -0.001281537845414626
$ws.Range("D29").Value = 0.02055877195361155
$ws.Range("E29").Value = -0.004379898570770036
$ws.Range("D30").Value = 0.02880823036340025
$ws.Range("E30").Value = 0.01791083801804172
$ws.Range("D31").Value = 0.01496990370390457
$ws.Range("E31").Value = 0.01357536265611659
$ws.Range("D32").Value = 0.01324680193114447
$ws.Range("E32").Value = -0.01260466372557845
$ws.Range("D33").Value = 0.01782502191971509
$ws.Range("E33").Value = 0.01331135902636937
$ws.Range("D34").Value = 0.04260541047867836
$ws.Range("E34").Value = 0.003991213143157646
$ws.Range("D35").Value = 0.01104648184806549
$ws.Range("E35").Value = 0.001363791339925147
$ws.Range("D36").Value = 0.01030038295474746
$ws.Range("E36").Value = -0.01160310553707022
$ws.Range("D37").Value = 0.01099987420644605
$ws.Range("E37").Value = -0.01262572223411074
$ws.Range("D38").Value = 0.007519366181269265
$ws.Range("E38").Value = -0.002754820936639257
$ws.Range("D39").Value = 0.01230592389715962
$ws.Range("E39").Value = -0.004506641366223807
$ws.Range("D40").Value = 0.01760915164425482
$ws.Range("E40").Value = -0.008690723597203909
$ws.Range("D41").Value = 0.01721356732440877
$ws.Range("E41").Value = -0.01087784937296499
$ws.Range("D42").Value = 0.03204816763233934
$ws.Range("E42").Value = 0.004788906190322351
$ws.Range("D43").Value = 0.01143396867866049
$ws.Range("E43").Value = -0.004380932512777824
$ws.Range("D44").Value = 0.02194692641905674
$ws.Range("E44").Value = -0.004344549363833883
$ws.Range("D45").Value = 0.01225187786425143
$ws.Range("E45").Value = 0.01839811868861552
$ws.Range("D46").Value = 0.0086450113440175
$ws.Range("E46").Value = 0.006142787126286464
$ws.Range("D47").Value = 0.01352864477406875
$ws.Range("E47").Value = -0.01085730989267975
$ws.Range("D48").Value = 0.01076551780242428
$ws.Range("E48").Value = -0.01755352644836283
$ws.Range("D49").Value = 0.01582098748693103
$ws.Range("E49").Value = -0.001928250053562564
$ws.Range("D50").Value = 0.008664407655479324
$ws.Range("E50").Value = -0.01539503664018715
$ws.Range("D51").Value = 0.01100103547428438
$ws.Range("E51").Value = -0.01008239375542053
$ws.Range("D52").Value = 0.008361065664737329
$ws.Range("E52").Value = 0.001899413659261606
$ws.Range("D53").Value = 0.009874605671104842
$ws.Range("E53").Value = -0.01081297557068472
$ws.Range("D54").Value = 0.1363588315109588
$ws.Range("E54").Value = -0.00009851246182634998
$ws.Range("D55").Value = 0.04368752379029787
$ws.Range("E55").Value = -0.002320468978993628
$ws.Range("E56").Value = -0.0009986903409633374

$ws.Protect()
